$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the source inlineStr cells) instead
# of auto-converting to a numeric value.
$textCells = @("D5", "D6", "D9", "D10", "D12", "D15", "D19", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D30", "D32", "D33", "D37", "D38", "D39", "D42", "D44", "D45", "D46")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "45.607.50"
$ws.Range("E2").Value = "  +6.75%  "
$ws.Range("D3").Value = "2.395.39"
$ws.Range("E3").Value = "  +3.87%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "114.15"
$ws.Range("D6").Value = "318.87"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "42.07"
$ws.Range("E10").Value = "  +5.92%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value = "8.72"
$ws.Range("E12").Value = "  +5.10%  "
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "15.88"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").Value = "2.758.54"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "2.396.17"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "45.542.37"
$ws.Range("E18").Value = "  +6.24%  "
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("D21").Value = "13.45"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "74.74"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "264.56"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +4.65%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Value = "7.67"
$ws.Range("E27").Value = "  +4.12%  "
$ws.Range("D28").Value = "11.32"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").Value = "39.51"
$ws.Range("E30").Value = "  +5.06%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").Value = "0.0967"
$ws.Range("E32").Value = "  +12.11%  "
$ws.Range("D33").Value = "172.60"
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  +6.57%  "
$ws.Range("D37").Value = "0.117"
$ws.Range("E37").Value = "  +4.94%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.09"
$ws.Range("E38").Value = "  +12.34%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "3.07"
$ws.Range("E39").Value = "  +8.28%  "
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("E41").Value = "  +11.84%  "
$ws.Range("D42").Value = "101.63"
$ws.Range("E42").Value = "  -5.95%  "
$ws.Range("E43").Value = "  +4.96%  "
$ws.Range("D44").Value = "13.55"
$ws.Range("E44").Value = "  +9.57%  "
$ws.Range("D45").Value = "72.17"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "88.07"
$ws.Range("E46").Value = "  +15.71%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("E49").Value = "  +10.78%  "
$ws.Range("E50").Value = "  +7.25%  "
$ws.Range("D51").Value = "1.664.88"
$ws.Range("E51").Value = "  -3.33%  "
